$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-08-18 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-19 Saturday", 2)
$d.Content.Find.Execute("12÷6=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "13÷5=2, 3", 2)
$d.Content.Find.Execute("60÷3=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷5=3, 4", 2)
$d.Content.Find.Execute("51÷2=25, 1", $true, $false, $false, $false, $false, $true, 1, $false, "14÷7=2, 0", 2)
$d.Content.Find.Execute("99÷7=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "87÷7=12, 3", 2)
$d.Content.Find.Execute("61÷6=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "41÷8=5, 1", 2)
$d.Content.Find.Execute("26÷3=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "75÷8=9, 3", 2)
$d.Content.Find.Execute("33÷7=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "23÷7=3, 2", 2)
$d.Content.Find.Execute("33÷9=3, 6", $true, $false, $false, $false, $false, $true, 1, $false, "91÷4=22, 3", 2)
$d.Content.Find.Execute("20÷5=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "33÷9=3, 6", 2)
$d.Content.Find.Execute("29÷7=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "14÷4=3, 2", 2)
$d.Content.Find.Execute("21÷2=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "11÷3=3, 2", 2)
$d.Content.Find.Execute("59÷3=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "88÷3=29, 1", 2)
$d.Content.Find.Execute("82÷8=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "59÷8=7, 3", 2)
$d.Content.Find.Execute("21÷9=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "90÷9=10, 0", 2)
$d.Content.Find.Execute("70÷9=7, 7", $true, $false, $false, $false, $false, $true, 1, $false, "58÷6=9, 4", 2)
$d.Content.Find.Execute("36÷3=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "96÷8=12, 0", 2)
$d.Content.Find.Execute("76÷5=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "36÷3=12, 0", 2)
$d.Content.Find.Execute("33÷3=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "79÷9=8, 7", 2)
$d.Content.Find.Execute("54÷9=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "65÷4=16, 1", 2)
$d.Content.Find.Execute("22÷6=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "60÷6=10, 0", 2)
$d.Content.Find.Execute("40÷2=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "73÷6=12, 1", 2)
$d.Content.Find.Execute("23÷5=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "75÷5=15, 0", 2)
$d.Content.Find.Execute("73÷4=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "47÷7=6, 5", 2)
$d.Content.Find.Execute("48÷6=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "43÷8=5, 3", 2)
$d.Content.Find.Execute("72÷2=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "66÷3=22, 0", 2)
